$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells we are about to write so that
# numeric-looking strings (e.g. "577.50", "0.555", "1.00") are kept
# as literal text instead of being coerced into numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.832.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.99%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.460.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.89%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.98'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.47%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.460.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.82%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.62%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.66'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.07%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.99%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.74%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.052.85'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.95%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.61'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.68%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.460.42'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.70%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.803.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.87%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.28'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.82%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.15%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.00'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.11%  '

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.555'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.19%  '

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.41'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.21%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.12%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.604.29'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000115'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.27%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -8.44%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.01%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.07%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.69%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.04%  '

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.09%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.59'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.01%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.63%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.03'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.58%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.09%  '

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '169.58'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.70%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.500.04'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.02%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0751'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.87%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.798'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.33%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.35'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.32%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.59%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.97%  '

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.30%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.597.51'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.24'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +10.40%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.81'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.52%  '

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.47%  '
